$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename Population_name related columns (values unchanged here aside from reindex)
$ws.Range("B1").Value = "Population_name"
$ws.Range("H1").Value = "Files_to_upload"
$ws.Range("I1").Value = "Expected_File_names"

# Row 2 (scenario1) - update test/population identifiers
$ws.Range("B2").Value = "NewImportLogic_3 - Test_Automation_3 - 1/13/2023"
$ws.Range("C2").Value = "NewImportLogic_3 - Test_Automation_3"
$ws.Range("D2").Value = "NewImportLogic_3 - Test_Automation_3_radio_button"
$ws.Range("E2").Value = "Clinical"
$ws.Range("F2").Value = "Clinical_radio_button"
$ws.Range("G2").Value = "\Testdata\Templates\SLRReport_SourceData\LIVEHTA_723_Testdata\ExpectedData_withFA13data.xlsx"
$ws.Range("H2").Value = "\Testdata\Templates\SLRReport_SourceData\LIVEHTA_723_Testdata\Master-Extraction-Template-Oncology_FA13data_FA19NR.xlsx"
$ws.Range("I2").Value = "Master-Extraction-Template-Oncology_FA13data_FA19NR.xlsx"
$ws.Range("J2").Value = "ExcelReport-NewImportLogic_3 - Test_Automation_3-Clinical-"

# Row 3
$ws.Range("J3").Value = "WordReport-NewImportLogic_3 - Test_Automation_3-Clinical-"

# Row 4 (scenario2) - update test/population identifiers
$ws.Range("B4").Value = "NewImportLogic_3 - Test_Automation_3 - 1/13/2023"
$ws.Range("C4").Value = "NewImportLogic_3 - Test_Automation_3"
$ws.Range("D4").Value = "NewImportLogic_3 - Test_Automation_3_radio_button"
$ws.Range("E4").Value = "Clinical"
$ws.Range("F4").Value = "Clinical_radio_button"
$ws.Range("G4").Value = "\Testdata\Templates\SLRReport_SourceData\LIVEHTA_723_Testdata\ExpectedData_withFA13_and_FA19data.xlsx"
$ws.Range("H4").Value = "\Testdata\Templates\SLRReport_SourceData\LIVEHTA_723_Testdata\Master-Extraction-Template-Oncology_FA13data_FA19data.xlsx"
$ws.Range("I4").Value = "Master-Extraction-Template-Oncology_FA13data_FA19data.xlsx"

# Update selection / scroll position to match the new view (J3 selected, H1 at top-left)
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J3").Select()

# Force pageSetup block to be written out (matches target diff adding <pageSetup .../>)
$ws.PageSetup.Orientation = 1
